{"js": "// Replace the date line and each \"a\u00f7b=c, d\" table-cell answer with its\n// updated value, per the commit diff. All \"old\" values are unique within\n// the document, so a direct search-and-replace keyed on the old text is\n// safe and unambiguous.\nconst replacements = [\n  [\"2025-09-21 Sunday\", \"2025-09-22 Monday\"],\n  [\"308\u00f78=38, 4\", \"126\u00f77=18, 0\"],\n  [\"853\u00f76=142, 1\", \"790\u00f74=197, 2\"],\n  [\"170\u00f76=28, 2\", \"753\u00f77=107, 4\"],\n  [\"989\u00f74=247, 1\", \"869\u00f77=124, 1\"],\n  [\"183\u00f79=20, 3\", \"113\u00f78=14, 1\"],\n  [\"958\u00f73=319, 1\", \"149\u00f73=49, 2\"],\n  [\"560\u00f76=93, 2\", \"209\u00f77=29, 6\"],\n  [\"628\u00f72=314, 0\", \"552\u00f74=138, 0\"],\n  [\"232\u00f73=77, 1\", \"264\u00f79=29, 3\"],\n  [\"162\u00f74=40, 2\", \"712\u00f74=178, 0\"],\n  [\"896\u00f73=298, 2\", \"903\u00f75=180, 3\"],\n  [\"288\u00f75=57, 3\", \"329\u00f72=164, 1\"],\n  [\"905\u00f77=129, 2\", \"503\u00f74=125, 3\"],\n  [\"527\u00f79=58, 5\", \"492\u00f76=82, 0\"],\n  [\"316\u00f77=45, 1\", \"632\u00f79=70, 2\"],\n  [\"616\u00f74=154, 0\", \"938\u00f75=187, 3\"],\n  [\"722\u00f76=120, 2\", \"386\u00f78=48, 2\"],\n  [\"545\u00f74=136, 1\", \"660\u00f74=165, 0\"],\n  [\"896\u00f77=128, 0\", \"974\u00f77=139, 1\"],\n  [\"360\u00f77=51, 3\", \"459\u00f74=114, 3\"],\n  [\"480\u00f78=60, 0\", \"675\u00f73=225, 0\"],\n  [\"591\u00f75=118, 1\", \"500\u00f77=71, 3\"],\n  [\"204\u00f74=51, 0\", \"954\u00f79=106, 0\"],\n  [\"364\u00f73=121, 1\", \"996\u00f79=110, 6\"],\n  [\"893\u00f78=111, 5\", \"392\u00f79=43, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"a\u00f7b=c, d\" table-cell answer with its\n# updated value, per the commit diff. All \"old\" values are unique within\n# the document, so Find/Replace keyed on the old text is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-21 Sunday\", \"2025-09-22 Monday\"),\n    @(\"308\u00f78=38, 4\", \"126\u00f77=18, 0\"),\n    @(\"853\u00f76=142, 1\", \"790\u00f74=197, 2\"),\n    @(\"170\u00f76=28, 2\", \"753\u00f77=107, 4\"),\n    @(\"989\u00f74=247, 1\", \"869\u00f77=124, 1\"),\n    @(\"183\u00f79=20, 3\", \"113\u00f78=14, 1\"),\n    @(\"958\u00f73=319, 1\", \"149\u00f73=49, 2\"),\n    @(\"560\u00f76=93, 2\", \"209\u00f77=29, 6\"),\n    @(\"628\u00f72=314, 0\", \"552\u00f74=138, 0\"),\n    @(\"232\u00f73=77, 1\", \"264\u00f79=29, 3\"),\n    @(\"162\u00f74=40, 2\", \"712\u00f74=178, 0\"),\n    @(\"896\u00f73=298, 2\", \"903\u00f75=180, 3\"),\n    @(\"288\u00f75=57, 3\", \"329\u00f72=164, 1\"),\n    @(\"905\u00f77=129, 2\", \"503\u00f74=125, 3\"),\n    @(\"527\u00f79=58, 5\", \"492\u00f76=82, 0\"),\n    @(\"316\u00f77=45, 1\", \"632\u00f79=70, 2\"),\n    @(\"616\u00f74=154, 0\", \"938\u00f75=187, 3\"),\n    @(\"722\u00f76=120, 2\", \"386\u00f78=48, 2\"),\n    @(\"545\u00f74=136, 1\", \"660\u00f74=165, 0\"),\n    @(\"896\u00f77=128, 0\", \"974\u00f77=139, 1\"),\n    @(\"360\u00f77=51, 3\", \"459\u00f74=114, 3\"),\n    @(\"480\u00f78=60, 0\", \"675\u00f73=225, 0\"),\n    @(\"591\u00f75=118, 1\", \"500\u00f77=71, 3\"),\n    @(\"204\u00f74=51, 0\", \"954\u00f79=106, 0\"),\n    @(\"364\u00f73=121, 1\", \"996\u00f79=110, 6\"),\n    @(\"893\u00f78=111, 5\", \"392\u00f79=43, 5\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
